$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    3 = @(1.505614041169197, 0.05231270169004087, 157.8057217802531, 0.4998867070740569, 0, 159.8635352301864)
    4 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    5 = @(0.001754667048134761, 0.05231270169004087, 3.082599426703578, 246.9852506941017, 0, 250.1219174895434)
    6 = @(0.06328177979961902, 0.05231270169004087, 3.082599426703578, 6.48142807727062, 1, 9.679621985463859)
    7 = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 3.811642989160245)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 2]
    }
}
